$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")
$ws.Activate()

# --- Question 3 table (rows 56-61): VLOOKUP array formulas spilling B:D ---
$ws.Range("B56:B61").FormulaArray = '=VLOOKUP($A$56:$A$61,A2:D52,4,FALSE)'
$ws.Range("C56:C61").FormulaArray = '=VLOOKUP($A$56:$A$61,A2:I52,9,FALSE)'
$ws.Range("D56:D61").FormulaArray = '=VLOOKUP($A$56:$A$61,A2:N52,14,FALSE)'

# --- Question 4 table (rows 65-70): XLOOKUP array formulas spilling B:D ---
$ws.Range("B65:B70").FormulaArray = '=XLOOKUP($A$65:$A$70,A2:A52,D2:D52,"N/A",0,1)'
$ws.Range("C65:C70").FormulaArray = '=XLOOKUP($A$65:$A$70, A2:A52,I2:I52, "N/A",0,1)'
$ws.Range("D65:D70").FormulaArray = '=XLOOKUP($A$65:$A$70,A2:A52,N2:N52,"N/A",0,1)'

$excel.Calculate()

# --- View state: scroll so row 61 is the top-left row, and select D66 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 61
$win.ScrollColumn = 1
$ws.Range("D66").Select()
